$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Row 6: CC_cal_temp value
$ws.Range("H6").Value = 20.11

# Row 7: CC_wl array (long string)
$ws.Range("H7").Value = '[187.73, 188.54, 189.36, 190.17, 190.99, 191.8, 192.62, 193.43, 194.25, 195.06, 195.88, 196.69, 197.51, 198.32, 199.14, 199.95, 200.77, 201.59, 202.4, 203.22, 204.03, 204.85, 205.67, 206.48, 207.3, 208.12, 208.93, 209.75, 210.57, 211.39, 212.2, 213.02, 213.84, 214.66, 215.47, 216.29, 217.11, 217.93, 218.75, 219.56, 220.38, 221.2, 222.02, 222.84, 223.65, 224.47, 225.29, 226.11, 226.93, 227.75, 228.57, 229.39, 230.2, 231.02, 231.84, 232.66, 233.48, 234.3, 235.12, 235.94, 236.76, 237.58, 238.4, 239.22, 240.03, 240.85, 241.67, 242.49, 243.31, 244.13, 244.95, 245.77, 246.59, 247.41, 248.23, 249.05, 249.87, 250.69, 251.51, 252.33, 253.15, 253.97, 254.79, 255.61, 256.43, 257.25, 258.06, 258.88, 259.7, 260.52, 261.34, 262.16, 262.98, 263.8, 264.62, 265.44, 266.26, 267.08, 267.9, 268.72, 269.54, 270.35, 271.17, 271.99, 272.81, 273.63, 274.45, 275.27, 276.09, 276.91, 277.72, 278.54, 279.36, 280.18, 281, 281.82, 282.64, 283.45, 284.27, 285.09, 285.91, 286.73, 287.54, 288.36, 289.18, 290, 290.81, 291.63, 292.45, 293.27, 294.08, 294.9, 295.72, 296.53, 297.35, 298.17, 298.98, 299.8, 300.62, 301.43, 302.25, 303.06, 303.88, 304.7, 305.51, 306.33, 307.14, 307.96, 308.77, 309.59, 310.4, 311.22, 312.03, 312.85, 313.66, 314.47, 315.29, 316.1, 316.91, 317.73, 318.54, 319.35, 320.17, 320.98, 321.79, 322.61, 323.42, 324.23, 325.04, 325.85, 326.67, 327.48, 328.29, 329.1, 329.91, 330.72, 331.53, 332.34, 333.15, 333.96, 334.77, 335.58, 336.39, 337.2, 338.01, 338.82, 339.63, 340.44, 341.25, 342.05, 342.86, 343.67, 344.48, 345.28, 346.09, 346.9, 347.7, 348.51, 349.32, 350.12, 350.93, 351.73, 352.54, 353.34, 354.15, 354.95, 355.76, 356.56, 357.37, 358.17, 358.97, 359.78, 360.58, 361.38, 362.18, 362.99, 363.79, 364.59, 365.39, 366.19, 366.99, 367.79, 368.59, 369.39, 370.19, 370.99, 371.79, 372.59, 373.39, 374.19, 374.98, 375.78, 376.58, 377.38, 378.17, 378.97, 379.77, 380.56, 381.36, 382.15, 382.95, 383.74, 384.54, 385.33, 386.12, 386.92, 387.71, 388.5, 389.3, 390.09, 390.88, 391.67, 392.46, 393.25, 394.05, 394.84]'

# Row 8: CC_eno3 array (long string)
$ws.Range("H8").Value = '[-0.00071178, -0.00070841, -0.00701958, 0.00416275, 0.00400766, 0.00565546, -0.00022933, 0.00371424, 0.00556231, 0.0051077, 0.00731, 0.00055434, 0.00028464, 0.00197316, 0.00197658, 0.00292466, -0.00019539, 0.00522575, 0.00242978, 0.0039731, 0.00123638, 0.00107095, 0.00078364, 0.00370094, 0.00565466, 0.00712946, 0.00674376, 0.00703045, 0.00693199, 0.00672029, 0.00650912, 0.00618648, 0.00589378, 0.00555232, 0.0052029, 0.00486638, 0.0045335, 0.00419383, 0.00385798, 0.00352514, 0.00319055, 0.00288161, 0.00257903, 0.00229822, 0.00204613, 0.00178655, 0.0015738, 0.00137019, 0.00118972, 0.00102497, 0.00086839, 0.00073153, 0.00062598, 0.00052109, 0.00042815, 0.00035533, 0.00029465, 0.00024071, 0.00019092, 0.00014172, 0.00012197, 0.00009214, 0.00006823, 0.00004905, 0.00003596, 0.0000262, 0.0000089, -0.00000592, -0.00000677, -0.00000954, -0.00001382, -0.0000109, -0.00001839, -0.00001475, -0.00000727, 0.00000184, 0.00000297, 0.00000226, -0.00001298, -0.00000472, -0.0000086, -0.00000366, -0.00000489, -0.00000629, 0.00000516, 0.00001088, 0.00000357, 0.00001356, 0.00001723, 0.00002946, 0.00002698, 0.00003762, 0.00002963, 0.00003793, 0.00003405, 0.00004386, 0.00004686, 0.00005426, 0.0000539, 0.00006828, 0.00005717, 0.00006891, 0.00006872, 0.00006826, 0.00006567, 0.00006587, 0.00008244, 0.00008321, 0.0000983, 0.00010253, 0.0001122, 0.00010295, 0.00010513, 0.00010821, 0.0001058, 0.00009387, 0.00010221, 0.00011223, 0.00011684, 0.00013476, 0.00013015, 0.0001387, 0.00014341, 0.00014531, 0.00015637, 0.00015633, 0.00016485, 0.00016009, 0.00016249, 0.00016204, 0.000173, 0.00016009, 0.00016921, 0.00017429, 0.00018078, 0.00018005, 0.00018419, 0.00019869, 0.00019913, 0.0002059, 0.00020866, 0.00021333, 0.00022009, 0.0002285, 0.00022142, 0.00022381, 0.0002255, 0.00021834, 0.00023277, 0.00023366, 0.00023178, 0.00024748, 0.00026255, 0.00026021, 0.00024963, 0.00026413, 0.00027171, 0.00024767, 0.00025339, 0.00025519, 0.00025447, 0.00025861, 0.00025519, 0.00027711, 0.00028176, 0.00028476, 0.00030805, 0.00029823, 0.00029476, 0.00029788, 0.00030489, 0.00031295, 0.00032694, 0.00032442, 0.00031405, 0.00032722, 0.00031589, 0.00034083, 0.0003257, 0.00034472, 0.00034134, 0.00033418, 0.00034343, 0.00034404, 0.00034647, 0.00036808, 0.0003734, 0.00036444, 0.00035636, 0.00036442, 0.00035248, 0.0003806, 0.00038387, 0.00039291, 0.00038647, 0.00038868, 0.00039475, 0.00040688, 0.00039993, 0.00039877, 0.00039753, 0.0003981, 0.00041601, 0.00041887, 0.00040798, 0.0004174, 0.00041151, 0.00042857, 0.00040042, 0.000427, 0.00043334, 0.00045976, 0.00044227, 0.00046136, 0.00052563, 0.00047205, 0.00046495, 0.00049565, 0.00043954, 0.00048607, 0.00041395, 0.00047635, 0.00043995, 0.00048355, 0.00049923, 0.00047738, 0.00046102, 0.00050587, 0.00056021, 0.00056691, 0.00055045, 0.00052834, 0.00053682, 0.00055073, 0.0005167, 0.00054407, 0.00051465, 0.00056171, 0.00055953, 0.00052307, 0.00062666, 0.00057828, 0.00059632, 0.00057622, 0.00059445, 0.00067006, 0.00063472, 0.0006565, 0.0006111, 0.00069391, 0.00066389, 0.00065506, 0.00067357, 0.00066778, 0.00067421, 0.00070101]'

# Row 9: CC_eswa array (long string)
$ws.Range("H9").Value = '[-0.00188787, 0.0103141, 0.01643845, 0.0116863, 0.00846114, 0.00556376, 0.01109636, 0.00429653, 0.00217128, 0.00660825, 0.00767925, 0.00893392, 0.01137808, 0.02204557, 0.03706705, 0.05028257, 0.0612044, 0.0654342, 0.07076886, 0.07746282, 0.0783656, 0.08083129, 0.07634903, 0.07167102, 0.06409565, 0.05487332, 0.04726339, 0.03972266, 0.03317938, 0.02749304, 0.02264431, 0.0185539, 0.01510482, 0.01221724, 0.00982781, 0.00785467, 0.0062264, 0.0049077, 0.00383849, 0.00298935, 0.00233602, 0.00181124, 0.00141481, 0.00110815, 0.00085262, 0.00068424, 0.0005389, 0.00043164, 0.00033956, 0.00027734, 0.00022651, 0.0001878, 0.00015918, 0.00012314, 0.00011276, 0.00009529, 0.00008682, 0.00007297, 0.00005444, 0.00004975, 0.00003683, 0.00003464, 0.00002695, 0.00001902, 0.00001045, 0.0000123, 0.00001552, 0.00000403, 0.0000088, 0.00000095, 0.00000479, -0.00000784, -0.0000074, -0.00002126, -0.00001426, -0.00002574, -0.00001964, -0.00000995, 0.00001027, 0.00000352, 0.00000252, -0.00000335, 0.00000082, 0.00000561, 0.00000413, -0.00000109, 0.00001466, 0.00000255, 0.00000959, 0.0000046, 0.00000807, 0.00001218, 0.00002077, 0.00002064, 0.00002604, 0.00002732, 0.00002468, 0.00003002, 0.00002531, 0.00002314, 0.00003226, 0.00003182, 0.0000496, 0.00005134, 0.00006118, 0.00004993, 0.00006627, 0.00006328, 0.00004723, 0.00006442, 0.00006612, 0.00005685, 0.00005825, 0.00004954, 0.00007881, 0.00007141, 0.00006859, 0.00006958, 0.00007994, 0.0000689, 0.00007937, 0.00007232, 0.00007724, 0.00008087, 0.00006746, 0.0000766, 0.00007798, 0.00007566, 0.00009279, 0.00009264, 0.00008873, 0.00009748, 0.00009366, 0.00010329, 0.00011158, 0.00010283, 0.00010559, 0.00010946, 0.00011042, 0.00011558, 0.00011438, 0.00013501, 0.00012083, 0.00012166, 0.00014446, 0.00014171, 0.00014762, 0.00016109, 0.00015604, 0.00015724, 0.00016292, 0.00015937, 0.00015189, 0.00017461, 0.00018931, 0.00017139, 0.00016665, 0.00018549, 0.0001985, 0.0002072, 0.00019868, 0.00022534, 0.00023477, 0.00022558, 0.00022159, 0.00022883, 0.00017794, 0.00021237, 0.00022534, 0.00021986, 0.00021997, 0.00023812, 0.00025031, 0.00024366, 0.00025788, 0.00024926, 0.00026095, 0.00025321, 0.00027143, 0.0002625, 0.0002735, 0.00027963, 0.00027655, 0.0002711, 0.00028352, 0.00027941, 0.00027875, 0.00029209, 0.00031962, 0.00031841, 0.00031174, 0.00030735, 0.00034145, 0.00033042, 0.00032408, 0.00033441, 0.00033784, 0.00035977, 0.00035182, 0.00037354, 0.00035175, 0.00037827, 0.00036641, 0.00036633, 0.00038846, 0.00039311, 0.00039217, 0.00042086, 0.00043134, 0.00041989, 0.0004121, 0.00042737, 0.0004487, 0.00039763, 0.00033855, 0.00042472, 0.0004204, 0.00041107, 0.00045284, 0.00045135, 0.00049394, 0.00044445, 0.00046376, 0.00044278, 0.00047579, 0.00045268, 0.00049773, 0.00051791, 0.00049441, 0.00051471, 0.00050354, 0.00049819, 0.00055551, 0.00052388, 0.00053371, 0.00055268, 0.00056897, 0.00057621, 0.00058677, 0.00058057, 0.00054196, 0.00057637, 0.00054678, 0.00059888, 0.00060846, 0.00059225, 0.00058051, 0.00067595, 0.00063375, 0.00063443, 0.00065827, 0.0007327, 0.00078675, 0.00082814, 0.0007928, 0.00084168]'

# Row 10: CC_di array (long string)
$ws.Range("H10").Value = '[21.75, 31.04166667, 38.375, 44.58333333, 50.125, 54.58333333, 55.5, 55.20833333, 55.375, 55.875, 53.45833333, 53.875, 65.54166667, 156.125, 527.08333333, 1597.54166667, 3855, 7362.91666667, 11519.75, 15495.45833333, 18779.08333333, 21377.20833333, 23519.08333333, 25516.5, 27598.66666667, 29924.33333333, 32565.41666667, 35470.08333333, 38672.45833333, 41876.41666667, 44899.33333333, 47537.91666667, 49366.125, 50225.04166667, 50015.41666667, 48800.79166667, 46821.625, 44333.25, 41753.45833333, 39219.08333333, 37015.875, 35111.625, 33590.625, 32440.66666667, 31753.91666667, 31381.16666667, 31425.33333333, 31843.5, 32583.875, 33645.45833333, 35066.83333333, 36735.08333333, 38691.375, 40838.08333333, 43066.125, 45261.66666667, 47354.41666667, 49113.58333333, 50407.5, 51093.75, 51073, 50403.875, 49070.375, 47223.29166667, 45056.54166667, 42593.95833333, 40114.29166667, 37713.70833333, 35461.70833333, 33368.04166667, 31568, 29960.58333333, 28585.29166667, 27435.08333333, 26470.45833333, 25703.66666667, 25099.5, 24666.54166667, 24372.41666667, 24252.16666667, 24238.625, 24374.25, 24622.20833333, 24963.45833333, 25403.875, 25914.125, 26464.95833333, 27037.54166667, 27592.875, 28086.20833333, 28474.41666667, 28777.54166667, 28873.375, 28788.16666667, 28472.20833333, 27929.625, 27183.41666667, 26273.29166667, 25180.83333333, 24011.04166667, 22842.58333333, 21663.125, 20523.41666667, 19442.5, 18445.125, 17546.125, 16743.75, 16043.45833333, 15452.70833333, 14951.29166667, 14540.70833333, 14235.58333333, 13979, 13806.45833333, 13701.79166667, 13665.45833333, 13683.79166667, 13761.875, 13897.58333333, 14094.08333333, 14359.95833333, 14670.875, 15047.29166667, 15480.125, 15968.20833333, 16487.5, 17043.79166667, 17634.625, 18239.75, 18845.125, 19427.04166667, 19991.33333333, 20501.08333333, 20943.08333333, 21317.08333333, 21609.95833333, 21802.79166667, 21916.70833333, 21906.66666667, 21792.08333333, 21614, 21339.29166667, 20997.45833333, 20581.75, 20119.5, 19617.83333333, 19105.125, 18550.58333333, 18013, 17470.83333333, 16928.20833333, 16427.875, 15940.5, 15475.95833333, 15058.70833333, 14653.29166667, 14272.75, 13925.45833333, 13619.875, 13306.625, 13023.33333333, 12751.29166667, 12499.20833333, 12277.95833333, 12064, 11868.91666667, 11705.04166667, 11555.20833333, 11431.29166667, 11321.5, 11233.20833333, 11153.54166667, 11100.25, 11053.625, 11044.41666667, 11034.79166667, 11037.41666667, 11052.20833333, 11080.16666667, 11109.625, 11150.58333333, 11200.75, 11238.625, 11275.16666667, 11298.83333333, 11322.08333333, 11357.95833333, 11375.75, 11373.33333333, 11363.04166667, 11344.66666667, 11310.91666667, 11257.125, 11204.25, 11136.33333333, 11055.5, 10947.625, 10823.04166667, 10680.95833333, 10527.66666667, 10352.625, 10168.54166667, 9958.45833333, 9739.375, 9507.79166667, 9259.41666667, 9012.29166667, 8745.54166667, 8481.58333333, 8211.91666667, 7942.5, 7669.54166667, 7408.16666667, 7156.54166667, 6899.625, 6659.20833333, 6433.79166667, 6214.375, 6009.95833333, 5825.41666667, 5641.5, 5471.29166667, 5299.375, 5127.95833333, 4966.625, 4817.5, 4682.58333333, 4567.16666667, 4459.54166667, 4365.29166667, 4290.91666667, 4223, 4161.33333333, 4116.79166667, 4070.45833333, 4034.54166667, 3996.33333333, 3945.95833333, 3890.95833333, 3834, 3788.04166667, 3744.70833333, 3711.625, 3677.875, 3650.04166667, 3628.70833333, 3610.125, 3598.29166667, 3593.70833333, 3588.25, 3582.29166667, 3580.125, 3580.45833333, 3585.95833333, 3581.70833333, 3527.75]'

# Row 11: CC_lower_wavelength_limit_for_spectra_fit value.
# Build the new "number, green fill" look by copying the format already used
# a few rows down (H16) for calibration-coefficient numbers, then tweak the
# number format / vertical alignment to match.
$ws.Range("H16").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H11").NumberFormat = "0.00"
$ws.Range("H11").VerticalAlignment = -4108
$ws.Range("H11").Value = 217

# Row 12: CC_upper_wavelength_limit_for_spectra_fit value - reuse the format
# we just finished building on H11 so both cells land on the same style.
$ws.Range("H11").Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H12").Value = 240

# Update the active selection on Asset_Cal_Info to match the saved view state
$ws.Range("H10").Select()
